$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New mini "input form" block (rows 19-21) plus a second entry further down
# (rows 35-36). Shared strings must be interned in the same order the user
# typed them so the shared-strings table lines up: toko, nama sales, barang,
# omz, agung, wafer2000c;20pcs, ec.

$ws.Range("C20").Value = "toko"
$ws.Range("C19").Value = "nama sales"
$ws.Range("D20").Value = "barang"
$ws.Range("E20").Value = "omz"
$ws.Range("C21").Value = "agung"
$ws.Range("D21").Value = "wafer2000c;20pcs"
$ws.Range("E21").Value = 200000
$ws.Range("D35").Value = "omz"
$ws.Range("D36").Value = "ec"
$ws.Range("E36").Value = 15

# Leave the selection where the user ended up after entering the last row.
[void]$ws.Range("E35").Select()
